$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 52263
$ws.Range("B2").Value = "Davi Lucca Almeida"
$ws.Range("C2").Value = "P&D"
$ws.Range("D2").Value = "Viagem de negócios"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 45082
$ws.Range("G2").Value = 3860.36

# Row 3
$ws.Range("A3").Value = 70366
$ws.Range("B3").Value = "Amanda Silva"
$ws.Range("C3").Value = "Operações"
$ws.Range("D3").Value = "Problemas pessoais"
$ws.Range("E3").Value = 7
$ws.Range("F3").Value = 45097
$ws.Range("G3").Value = 6129.06

# Row 4
$ws.Range("A4").Value = 67991
$ws.Range("B4").Value = "Ian Peixoto"
$ws.Range("C4").Value = "Financeiro"
$ws.Range("D4").Value = "Consulta médica"
$ws.Range("E4").Value = 6
$ws.Range("F4").Value = 45081
$ws.Range("G4").Value = 10469.68

# Row 5
$ws.Range("A5").Value = 86743
$ws.Range("B5").Value = "Laura Cardoso"
$ws.Range("C5").Value = "Engenharia"
$ws.Range("F5").Value = 45097
$ws.Range("G5").Value = 11436.88

# Row 6
$ws.Range("A6").Value = 84932
$ws.Range("B6").Value = "Sr. Enrico Porto"
$ws.Range("C6").Value = "P&D"
$ws.Range("D6").Value = "Problemas pessoais"
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 45096
$ws.Range("G6").Value = 9200.530000000001

# Row 7
$ws.Range("A7").Value = 80652
$ws.Range("B7").Value = "Clarice Santos"
$ws.Range("C7").Value = "Financeiro"
$ws.Range("D7").Value = "Doença"
$ws.Range("E7").Value = 4
$ws.Range("F7").Value = 45084
$ws.Range("G7").Value = 3835.08

# Row 8
$ws.Range("A8").Value = 51237
$ws.Range("B8").Value = "Paulo Moraes"
$ws.Range("C8").Value = "Engenharia"
$ws.Range("D8").Value = "Consulta médica"
$ws.Range("E8").Value = 4
$ws.Range("F8").Value = 45089
$ws.Range("G8").Value = 8656.030000000001

# Row 9
$ws.Range("A9").Value = 80806
$ws.Range("B9").Value = "Henrique da Rocha"
$ws.Range("C9").Value = "Engenharia"
$ws.Range("D9").Value = "Problemas pessoais"
$ws.Range("E9").Value = 2
$ws.Range("F9").Value = 45091
$ws.Range("G9").Value = 9084.110000000001

# Row 10
$ws.Range("A10").Value = 32143
$ws.Range("B10").Value = "Sra. Maria Vitória Martins"
$ws.Range("D10").Value = "Outros"
$ws.Range("E10").Value = 4
$ws.Range("F10").Value = 45093
$ws.Range("G10").Value = 3605.29

# Row 11
$ws.Range("A11").Value = 7325
$ws.Range("B11").Value = "Lucas Mendes"
$ws.Range("C11").Value = "Marketing"
$ws.Range("D11").Value = "Viagem de negócios"
$ws.Range("E11").Value = 6
$ws.Range("F11").Value = 45084
$ws.Range("G11").Value = 3377.52
